# Update cryptos list: refresh price/volume figures and fix two swapped
# row pairs (Algorand/LidoDAOToken and TheGraph/RocketPoolETH) that were
# in the wrong rank order.
#
# Price cells (column D) are stored as text in the source data (e.g.
# "43.610.14" or "269.27"), so a leading apostrophe is used to force
# the Value assignment to stay text instead of being auto-converted to
# a number by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.652.00"
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = "'2.324.37"
$ws.Range("E3").Value = '  +4.27%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'270.67"
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = "'95.18"
$ws.Range("E6").Value = '  +6.62%  '
$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = '  +1.67%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = "'0.622"
$ws.Range("E9").Value = '  +2.57%  '
$ws.Range("D10").Value = "'45.26"
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("D11").Value = "'0.0945"
$ws.Range("E11").Value = '  +2.62%  '
$ws.Range("D12").Value = "'8.11"
$ws.Range("E12").Value = '  +2.34%  '
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").Value = "'2.670.19"
$ws.Range("E14").Value = '  +3.97%  '
$ws.Range("D15").Value = "'15.73"
$ws.Range("E15").Value = '  +4.35%  '
$ws.Range("D16").Value = "'0.863"
$ws.Range("E16").Value = '  +9.09%  '
$ws.Range("D17").Value = "'2.320.50"
$ws.Range("E17").Value = '  +3.96%  '
$ws.Range("D18").Value = "'43.617.47"
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").Value = "'0.0000109"
$ws.Range("E19").Value = '  +5.31%  '
$ws.Range("D20").Value = "'6.47"
$ws.Range("E20").Value = '  +8.00%  '
$ws.Range("D21").Value = "'72.05"
$ws.Range("E21").Value = '  +2.52%  '
$ws.Range("D22").Value = "'240.64"
$ws.Range("E22").Value = '  +3.69%  '
$ws.Range("D23").Value = "'2.27"
$ws.Range("E23").Value = '  -3.47%  '
$ws.Range("D24").Value = "'9.42"
$ws.Range("E24").Value = '  +8.60%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Value = "'11.44"
$ws.Range("E26").Value = '  +4.04%  '
$ws.Range("D27").Value = "'2.53"
$ws.Range("E27").Value = '  +1.97%  '
$ws.Range("D28").Value = "'3.47"
$ws.Range("E28").Value = '  -1.67%  '
$ws.Range("D30").Value = "'38.44"
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").Value = "'22.65"
$ws.Range("E31").Value = '  +9.36%  '
$ws.Range("D32").Value = "'172.73"
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("D33").Value = "'0.0905"
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("D34").Value = "'5.50"
$ws.Range("E34").Value = '  +2.51%  '
$ws.Range("E35").Value = '  +2.68%  '
$ws.Range("D36").Value = "'0.0361"
$ws.Range("E36").Value = '  +2.97%  '
$ws.Range("D37").Value = "'0.108"
$ws.Range("E37").Value = '  -2.27%  '
$ws.Range("D38").Value = "'4.41"
$ws.Range("E38").Value = '  +3.67%  '
$ws.Range("D39").Value = "'3.38"
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = "'0.237"
$ws.Range("E40").Value = '  +12.11%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").Value = "'2.34"
$ws.Range("E41").Value = '  +8.75%  '
$ws.Range("D42").Value = "'1.36"
$ws.Range("E42").Value = '  +18.99%  '
$ws.Range("D43").Value = "'12.10"
$ws.Range("E43").Value = '  -3.39%  '
$ws.Range("E44").Value = '  +7.74%  '
$ws.Range("D45").Value = "'62.20"
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("D46").Value = "'5.35"
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("E47").Value = '  +4.62%  '
$ws.Range("D48").Value = "'100.67"
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").Value = "'0.189"
$ws.Range("E50").Value = '  +17.12%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = "'2.550.20"
$ws.Range("E51").Value = '  +3.91%  '
